# Updated Documentation on RTOS Threads
# -------------------------------------------------
# Adds a "Thread / Type / Func / Description" reference table (rows 9-21)
# to the existing Events/Clock/HWI/SWI/Task summary sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate new documentation rows (9-21) in the precise order needed so the
# shared-strings table is built up in the same sequence as the target workbook.
$ws.Range('B9').Value2 = 'Type '
$ws.Range('C9').Value2 = 'Func'
$ws.Range('D9').Value2 = 'Description'
$ws.Range('A9').Value2 = 'Thread'
$ws.Range('A10').Value2 = 'SPI_Event'
$ws.Range('B10').Value2 = 'Event'
$ws.Range('C10').Value2 = 'Event_post'
$ws.Range('A13').Value2 = 'I2C_Event'
$ws.Range('B13').Value2 = 'Event'
$ws.Range('C13').Value2 = 'Event_post'
$ws.Range('D13').Value2 = 'Recieves output of posted I2C events '
$ws.Range('D10').Value2 = 'Recieves output of posted SPI events '
$ws.Range('A11').Value2 = 'I2C_Recieve_Event'
$ws.Range('B11').Value2 = 'Event'
$ws.Range('C11').Value2 = 'Event_post'
$ws.Range('D11').Value2 = 'Recieves output of posted I2C receive events '
$ws.Range('A12').Value2 = 'StateChangeEvent'
$ws.Range('B12').Value2 = 'Event'
$ws.Range('C12').Value2 = 'Event_post'
$ws.Range('D12').Value2 = 'Recieves output of posted StateChange events '
$ws.Range('A14').Value2 = 'CAN_RTR'
$ws.Range('B14').Value2 = 'Clock'
$ws.Range('C14').Value2 = 'SendCAN'
$ws.Range('D14').Value2 = 'Repeatadly calls sendCAN with a period of 500millisec and a timeoutof 1000millisec'
$ws.Range('A15').Value2 = 'SPI_HWI'
$ws.Range('B15').Value2 = 'HWI'
$ws.Range('C15').Value2 = 'SPI_HandleInterrupt'
$ws.Range('A16').Value2 = 'Timer_HWI'
$ws.Range('B16').Value2 = 'HWI'
$ws.Range('C16').Value2 = 'Timer_ISR '
$ws.Range('D16').Value2 = 'Increments timer'
$ws.Range('A17').Value2 = 'TCA9555_INT_HWI'
$ws.Range('B17').Value2 = 'HWI'
$ws.Range('C17').Value2 = 'I2C_TCA9555Interupt'
$ws.Range('A18').Value2 = 'I2C_TCA9555_HWI'
$ws.Range('B18').Value2 = 'HWI'
$ws.Range('C18').Value2 = 'I2C_Interrupt'
$ws.Range('D17').Value2 = 'Sets state for new inputs from TCA9555 to True'
$ws.Range('D18').Value2 = 'Interrupt c2000 generates for when you are allowed to modify registers for I2C'
$ws.Range('A19').Value2 = 'CAN_Recieve_HWI'
$ws.Range('B19').Value2 = 'HWI'
$ws.Range('C19').Value2 = 'CAN_Recieve_Interrupt'
$ws.Range('D19').Value2 = 'Checks for timeout /mailbox needs data to be read out of it | disables the mailbox once read from it  | Swaps mailbox message ID and re-enables it '
$ws.Range('A20').Value2 = 'SPI_EventHandleTask'
$ws.Range('B20').Value2 = 'Task'
$ws.Range('C20').Value2 = 'SPI_HandleEvent'
$ws.Range('D20').Value2 = 'Processes pending SPI event, which is a transmition ready event and an spi done event'
$ws.Range('A21').Value2 = 'BatteryController'
$ws.Range('B21').Value2 = 'Task'
$ws.Range('C21').Value2 = 'BatteryController_Task'
$ws.Range('D21').Value2 = 'Task run forever | triggers I2c update | fetches state | checks THEN processes states | '

# Header row (row 9) is bold, matching the existing header row 1 style.
$ws.Range('A9:D9').Font.Bold = $true

# The "Event_post" entry in C10 carries the "Calibri (Body)" font that shows
# up when the author pasted it in from another document.
$ws.Range('C10').Font.Name = 'Calibri (Body)'

# Match the author's final selection/cursor position.
$ws.Range('D21').Select() | Out-Null
